$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $ws.Range($cellRef).Value = "'" + $val
    $ws.Range($cellRef).Style = "Normal"
}

Set-TextCell "D2" "26.994.00"
Set-TextCell "E2" "  -0.71%  "
Set-TextCell "D3" "1.829.46"
Set-TextCell "E3" "  +0.06%  "
Set-TextCell "E4" "  -0.35%  "
Set-TextCell "D5" "311.64"
Set-TextCell "E5" "  -0.61%  "
Set-TextCell "E6" "  -0.32%  "
Set-TextCell "D7" "0.4656"
Set-TextCell "E7" "  -1.35%  "
Set-TextCell "D8" "0.3712"
Set-TextCell "E8" "  +1.40%  "
Set-TextCell "D9" "0.07447"
Set-TextCell "E9" "  +0.37%  "
Set-TextCell "D10" "0.8756"
Set-TextCell "E10" "  -0.78%  "
Set-TextCell "D11" "20.04"
Set-TextCell "E11" "  -1.93%  "
Set-TextCell "D12" "0.07872"
Set-TextCell "E12" "  +7.54%  "
Set-TextCell "D13" "1.833.13"
Set-TextCell "E13" "  -4.64%  "
Set-TextCell "D14" "5.368"
Set-TextCell "E14" "  -0.69%  "
Set-TextCell "E15" "  +1.01%  "
Set-TextCell "D16" "92.08"
Set-TextCell "E16" "  -1.49%  "
Set-TextCell "E17" "  +0.02%  "
Set-TextCell "D18" "0.000008943"
Set-TextCell "E18" "  +1.89%  "
Set-TextCell "E19" "  -0.38%  "
Set-TextCell "D20" "14.74"
Set-TextCell "E20" "  +0.17%  "
Set-TextCell "D21" "27.031.62"
Set-TextCell "E21" "  -2.30%  "
Set-TextCell "D22" "5.171"
Set-TextCell "E22" "  -1.90%  "
Set-TextCell "E23" "  -0.14%  "
Set-TextCell "D24" "2.062.93"
Set-TextCell "E24" "  -2.58%  "
Set-TextCell "D25" "152.75"
Set-TextCell "E25" "  +0.68%  "
Set-TextCell "D26" "1.831"
Set-TextCell "E26" "  -2.99%  "
Set-TextCell "E27" "  -1.69%  "
Set-TextCell "D28" "2.103"
Set-TextCell "E28" "  -1.68%  "
Set-TextCell "D29" "5.135"
Set-TextCell "E29" "  -1.23%  "
Set-TextCell "D30" "115.74"
Set-TextCell "E30" "  -0.83%  "
Set-TextCell "D31" "0.08883"
Set-TextCell "E31" "  -0.85%  "
Set-TextCell "D32" "2.969"
Set-TextCell "E32" "  +0.75%  "
Set-TextCell "D33" "0.7300"
Set-TextCell "E33" "  -1.87%  "
Set-TextCell "E34" "  -1.49%  "
Set-TextCell "E35" "  -2.84%  "
Set-TextCell "D36" "2.508"
Set-TextCell "E36" "  +3.95%  "
Set-TextCell "E37" "  -1.01%  "
Set-TextCell "D38" "0.01958"
Set-TextCell "E38" "  +0.46%  "
Set-TextCell "D39" "0.05251"
Set-TextCell "E39" "  -1.31%  "
Set-TextCell "D40" "7.356"
Set-TextCell "E40" "  +2.11%  "
Set-TextCell "D41" "2.924"
Set-TextCell "E41" "  -0.89%  "
Set-TextCell "D42" "0.5206"
Set-TextCell "E42" "  -1.35%  "
Set-TextCell "D43" "0.1629"
Set-TextCell "E43" "  -1.22%  "
Set-TextCell "D44" "0.8572"
Set-TextCell "E44" "  -15.22%  "
Set-TextCell "D45" "8.245"
Set-TextCell "E45" "  -2.42%  "
Set-TextCell "D46" "0.4864"
Set-TextCell "E46" "  -0.69%  "
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextCell "D47" "1.008"
Set-TextCell "E47" "  -0.31%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell "D48" "10.18"
Set-TextCell "E48" "  -2.54%  "
Set-TextCell "D49" "102.75"
Set-TextCell "E49" "  -2.06%  "
Set-TextCell "D50" "1.624"
Set-TextCell "E50" "  -1.83%  "
Set-TextCell "D51" "0.06246"
Set-TextCell "E51" "  -0.86%  "
